# Apply the periodic "cryptos list" refresh (GitHub Actions style update).
# Column A (rank) and all row/column structure stay untouched; only the
# Coin/Link labels (B/C) for two re-ranked coin pairs and the Price/Volume(1h)
# text values (D/E) are refreshed to their newly scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link cells that changed rank position -------------------------
$textUpdates = [ordered]@{
    "B32" = "Fetch.AI"
    "C32" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "B33" = "Aptos"
    "C33" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "B50" = "Arweave"
    "C50" = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
    "B51" = "FirstDigitalUSD"
    "C51" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
}

foreach ($coord in $textUpdates.Keys) {
    $ws.Range($coord).Value = $textUpdates[$coord]
}

# --- Price (D) / Volume(1h) (E) cells ---------------------------------------
# These are stored as plain text in the workbook (e.g. "68.127.75", "0.999",
# "  +0.00%  "), so force a text number format before assigning the value,
# otherwise Excel would silently reinterpret numeric-looking strings as
# numbers (dropping formatting such as trailing zeros). Reset the style back
# to "Normal" afterwards so no cell formatting is left behind.
$priceVolumeUpdates = [ordered]@{
    "D2" = "68.127.75"
    "E2" = "  +0.00%  "
    "D3" = "3.269.58"
    "E3" = "  +0.39%  "
    "D5" = "584.36"
    "E5" = "  +0.26%  "
    "D6" = "184.41"
    "E6" = "  -0.40%  "
    "E7" = "  +0.01%  "
    "E8" = "  +0.24%  "
    "E9" = "  -2.05%  "
    "D10" = "6.64"
    "E10" = "  -0.35%  "
    "E11" = "  -2.84%  "
    "D12" = "3.839.10"
    "E12" = "  +0.27%  "
    "D14" = "68.117.35"
    "E14" = "  -0.17%  "
    "D15" = "27.34"
    "E15" = "  -2.93%  "
    "E16" = "  -1.89%  "
    "D17" = "3.279.17"
    "E17" = "  +0.68%  "
    "E18" = "  -2.30%  "
    "D19" = "13.28"
    "E19" = "  -2.48%  "
    "D20" = "416.08"
    "E20" = "  +5.74%  "
    "D21" = "7.52"
    "E22" = "  +0.30%  "
    "D23" = "71.12"
    "E23" = "  -0.31%  "
    "E24" = "  -2.24%  "
    "D25" = "0.0000117"
    "E25" = "  -2.20%  "
    "E26" = "  -0.81%  "
    "E27" = "  -4.28%  "
    "D28" = "0.999"
    "E28" = "  -0.11%  "
    "D29" = "1.94"
    "E29" = "  -1.80%  "
    "D30" = "22.63"
    "E30" = "  -1.25%  "
    "E31" = "  -4.90%  "
    "D32" = "1.24"
    "E32" = "  -3.16%  "
    "D33" = "6.83"
    "E33" = "  -4.51%  "
    "D34" = "164.46"
    "E34" = "  +1.28%  "
    "E35" = "  -4.43%  "
    "E36" = "  -3.91%  "
    "D37" = "26.64"
    "E37" = "  -1.13%  "
    "D38" = "0.791"
    "E38" = "  -3.86%  "
    "D39" = "4.44"
    "E39" = "  -3.43%  "
    "E40" = "  -4.17%  "
    "D41" = "2.629.55"
    "E41" = "  -0.94%  "
    "E42" = "  -2.27%  "
    "E43" = "  -3.70%  "
    "D44" = "334.86"
    "E44" = "  -1.54%  "
    "D45" = "24.20"
    "E45" = "  -4.72%  "
    "D46" = "0.0274"
    "E46" = "  -3.05%  "
    "E47" = "  -0.45%  "
    "D48" = "6.22"
    "E48" = "  -2.18%  "
    "E49" = "  -1.62%  "
    "D50" = "30.70"
    "E50" = "  -2.74%  "
    "D51" = "1.00"
    "E51" = "  +0.00%  "
}

foreach ($coord in $priceVolumeUpdates.Keys) {
    $cell = $ws.Range($coord)
    $cell.NumberFormat = "@"
    $cell.Value = $priceVolumeUpdates[$coord]
    $cell.Style = "Normal"
}
